# Update "想去人数" (F column) counts across sheets to reflect newly
# generated data (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4693
$ws1.Range("F3").Value = 2549
$ws1.Range("F4").Value = 66
$ws1.Range("F6").Value = 239
$ws1.Range("F7").Value = 140
$ws1.Range("F8").Value = 211
$ws1.Range("F9").Value = 181
$ws1.Range("F10").Value = 1784
$ws1.Range("F11").Value = 327
$ws1.Range("F12").Value = 4079
$ws1.Range("F13").Value = 48
$ws1.Range("F14").Value = 281

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 47

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4693
$ws4.Range("F3").Value = 2549
$ws4.Range("F5").Value = 66
$ws4.Range("F7").Value = 47
$ws4.Range("F8").Value = 239
$ws4.Range("F9").Value = 140
$ws4.Range("F10").Value = 211
$ws4.Range("F11").Value = 181
$ws4.Range("F14").Value = 1784
$ws4.Range("F15").Value = 327
$ws4.Range("F16").Value = 4079
$ws4.Range("F17").Value = 48
$ws4.Range("F18").Value = 281
